$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("seats")

# Update existing rows 1-3 with new values
$ws.Cells.Item(1, 1).Value = 6
$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(3, 1).Value = 17

# Insert four new rows before the old row 4 ("23"/"admin"), pushing it down to row 8
$ws.Rows.Item(4).Resize(4).Insert()

# Fill in the newly inserted rows with the seat data
$ws.Cells.Item(4, 1).Value = 18
$ws.Cells.Item(4, 2).Value = "admin"

$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = "klf"

$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "klf"

$ws.Cells.Item(7, 1).Value = 25
$ws.Cells.Item(7, 2).Value = "admin"
